$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove AgTests (F) / AgPosit (G) values for rows that no longer have this data
$removeRows = @(20, 118, 178, 190, 211, 221, 222, 223, 224, 225, 226, 227, 228, 229, 230, 231, 232, 233, 234, 235, 236, 237, 238, 239, 240, 241, 242, 243, 244, 245, 246, 247, 248, 249, 250, 251, 252, 253, 254, 255, 256, 257, 258, 259, 260, 261, 262, 263, 264, 265, 266, 267, 268, 269, 270, 271, 272, 273, 274, 275, 276, 277, 278, 279, 280, 281, 282, 283, 284, 285, 286, 287, 288, 289, 290, 291, 292, 293, 294, 295, 296, 297, 298, 299, 300, 301, 302)
foreach ($r in $removeRows) {
    $ws.Cells.Item($r, 6).ClearContents()
    $ws.Cells.Item($r, 7).ClearContents()
}

# Add AgTests (F) / AgPosit (G) values for newly reported rows
$newData = @{
    849 = @(2045, 164)
    850 = @(642, 49)
    851 = @(992, 56)
    852 = @(2959, 272)
    853 = @(898, 92)
    854 = @(2202, 294)
    855 = @(2452, 201)
    856 = @(1886, 228)
    857 = @(657, 80)
    858 = @(954, 112)
    859 = @(3401, 439)
    860 = @(2459, 293)
    861 = @(2253, 290)
    862 = @(2687, 249)
    863 = @(2626, 331)
    864 = @(892, 118)
    865 = @(1006, 158)
    866 = @(4017, 617)
    867 = @(2802, 368)
    868 = @(2817, 332)
    869 = @(3348, 365)
    870 = @(3056, 384)
    871 = @(988, 142)
    872 = @(1194, 152)
    873 = @(5970, 694)
    874 = @(4406, 435)
    875 = @(2980, 331)
    876 = @(3924, 413)
    877 = @(3075, 388)
    878 = @(989, 103)
    879 = @(1303, 176)
    880 = @(4891, 568)
    881 = @(3860, 366)
    882 = @(3468, 349)
    883 = @(3672, 327)
    884 = @(2453, 257)
    885 = @(792, 59)
    886 = @(562, 59)
}
foreach ($r in $newData.Keys) {
    $pair = $newData[$r]
    $ws.Cells.Item($r, 6).Value = $pair[0]
    $ws.Cells.Item($r, 7).Value = $pair[1]
}

Write-Host "Done"